$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain text so numeric-looking values
# (e.g. "606.08") are not auto-converted to numbers by Excel, matching
# the inline-string storage used by the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "66.124.67"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "3.555.92"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "606.08"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").Value = "144.48"
$ws.Range("D7").Value = "3.555.94"
$ws.Range("E7").Value = "  +3.98%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  +4.55%  "
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("D11").Value = "7.91"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "4.158.05"
$ws.Range("E13").Value = "  +4.06%  "
$ws.Range("D15").Value = "30.06"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "3.551.92"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").Value = "66.227.09"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "11.30"
$ws.Range("E19").Value = "  +9.73%  "
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "14.82"
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("D22").Value = "429.96"
$ws.Range("E22").Value = "  +3.54%  "
$ws.Range("D23").Value = "0.611"
$ws.Range("E23").Value = "  +6.01%  "
$ws.Range("D24").Value = "79.20"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").Value = "3.696.64"
$ws.Range("E25").Value = "  +4.19%  "
$ws.Range("E27").Value = "  +7.76%  "
$ws.Range("E28").Value = "  +3.98%  "
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "9.07"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "1.46"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "25.51"
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("D34").Value = "3.549.37"
$ws.Range("E34").Value = "  +3.95%  "
$ws.Range("E35").Value = "  -5.44%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "1.74"
$ws.Range("E37").Value = "  +3.80%  "
$ws.Range("D38").Value = "7.88"
$ws.Range("E38").Value = "  +4.80%  "
$ws.Range("D39").Value = "5.61"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "174.90"
$ws.Range("E41").Value = "  +3.73%  "
$ws.Range("D42").Value = "0.0850"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "5.21"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").Value = "0.893"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("D45").Value = "1.92"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").Value = "46.05"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").Value = "25.80"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "23.70"
$ws.Range("E49").Value = "  +16.27%  "
$ws.Range("D50").Value = "7.12"
$ws.Range("E51").Value = "  +3.03%  "
